$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, copying the header formatting (bold,
# bordered, centered style) from the neighboring "sum" header in G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save flag values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
